# Updated the Web Inspect
# Adds two new locator sheets (ImageLocators, HeadingLocators) after the
# existing TableLocators sheet, each seeded with the same 7-column header
# row ("Loc1".."Loc7") used by the other *Locators sheets, and updates the
# selection/active-tab state accordingly.

$wb = $excel.ActiveWorkbook

# The template sheet whose header row / layout is reused for the new sheets.
$tableSheet = $wb.Worksheets.Item("TableLocators")

# --- Create "ImageLocators" right after "TableLocators" -------------------
$tableSheet.Copy([System.Type]::Missing, $tableSheet)
$imageSheet = $wb.Worksheets.Item($tableSheet.Index + 1)
$imageSheet.Name = "ImageLocators"

# --- Create "HeadingLocators" right after "ImageLocators" ------------------
$tableSheet.Copy([System.Type]::Missing, $imageSheet)
$headingSheet = $wb.Worksheets.Item($imageSheet.Index + 1)
$headingSheet.Name = "HeadingLocators"

# --- Restore TableLocators' own selection (it is no longer the active tab) -
$tableSheet.Activate() | Out-Null
$tableSheet.Range("A1:G1").Select() | Out-Null

# --- ImageLocators: header row selected, not the active tab ---------------
$imageSheet.Activate() | Out-Null
$imageSheet.Range("A1:G1").Select() | Out-Null

# --- HeadingLocators: becomes the active tab, cell B2 selected -------------
$headingSheet.Activate() | Out-Null
$headingSheet.Range("B2").Select() | Out-Null
